$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = 1.248583466539117
    "C2" = 0.5135445598622255
    "D2" = 0.02697253705651548
    "E2" = 0.4183939171993387
    "F2" = 1.40831938304018
    "I2" = 0.8201472520215134
    "B3" = 1.109004199223932
    "C3" = 0.4494248925602164
    "D3" = 0.02650423200020668
    "E3" = 0.3645020341276677
    "F3" = 1.3348920828372
    "I3" = 0.793131367362335
    "B4" = 1.023753385397583
    "C4" = 0.4101930759799188
    "D4" = 0.0262425623976732
    "E4" = 0.3315723959174477
    "F4" = 1.290923878185708
    "I4" = 0.7771895207634003
    "B5" = 0.9891235662772715
    "C5" = 0.3942382203240413
    "D5" = 0.02614226112430273
    "E5" = 0.3181899391099563
    "F5" = 1.273281674887869
    "I5" = 0.7708526310546588
    "B6" = 0.9833799061866557
    "C6" = 0.3915908223520432
    "D6" = 0.02612598349378459
    "E6" = 0.3159698997834681
    "F6" = 1.270368660628677
    "I6" = 0.7698099468155704
    "B7" = 1.023285911619382
    "C7" = 0.4099777748626821
    "D7" = 0.02624118429420008
    "E7" = 0.3313917718563317
    "F7" = 1.290684842283625
    "I7" = 0.7771034170905438
    "B8" = 1.200361005336902
    "C8" = 0.4914060859751999
    "D8" = 0.0268055902312696
    "E8" = 0.3997765224959551
    "F8" = 1.382766553573703
    "I8" = 0.8106963311955866
    "B9" = 1.551340953639112
    "C9" = 0.65229212190917
    "D9" = 0.02812564034314846
    "E9" = 0.5353223072584257
    "F9" = 1.572458603900571
    "I9" = 0.8818331445848173
    "B10" = 1.811729519971095
    "C10" = 0.7713974297906248
    "D10" = 0.02923668437058069
    "E10" = 0.6360383629342294
    "F10" = 1.717764074275038
    "I10" = 0.9374986201205786
    "B11" = 1.930791820376669
    "C11" = 0.8258153070449907
    "D11" = 0.02977525142903659
    "E11" = 0.6821581709987896
    "F11" = 1.785240746264435
    "I11" = 0.9636035652325461
    "B12" = 1.975969605603495
    "C12" = 0.8464588475313803
    "D12" = 0.02998416302251883
    "E12" = 0.6996706982716745
    "F12" = 1.810997008541079
    "I12" = 0.9736047686078706
    "B13" = 1.966235631998245
    "C13" = 0.8420112170801985
    "D13" = 0.02993894597989311
    "E13" = 0.6958968596913166
    "F13" = 1.805440747184122
    "I13" = 0.9714456257839004
    "B14" = 1.934506768806273
    "C14" = 0.8275129128982712
    "D14" = 0.02979233798153302
    "E14" = 0.6835979510093608
    "F14" = 1.787355597873272
    "I14" = 0.9644240297522373
    "B15" = 1.915083949766824
    "C15" = 0.8186371398834922
    "D15" = 0.02970318935458494
    "E15" = 0.6760708892048086
    "F15" = 1.776304720054412
    "I15" = 0.9601382827241025
    "B16" = 1.803961113385128
    "C16" = 0.7678460767689899
    "D16" = 0.02920217189517871
    "E16" = 0.633030799230653
    "F16" = 1.713382555926671
    "I16" = 0.9358086273887807
    "B17" = 1.735949755135266
    "C17" = 0.7367498156958163
    "D17" = 0.02890344601022576
    "E17" = 0.6067078802236949
    "F17" = 1.675138899022784
    "I17" = 0.921085915907824
    "B18" = 1.696888627717556
    "C18" = 0.7188861426147923
    "D18" = 0.02873473454552311
    "E18" = 0.591595943341261
    "F18" = 1.653271332854302
    "I18" = 0.9126912477417193
    "B19" = 1.683672915268005
    "C19" = 0.7128415201045186
    "D19" = 0.02867813942287967
    "E19" = 0.5864840339517627
    "F19" = 1.645889343899967
    "I19" = 0.9098614716022126
    "B20" = 1.743183733073693
    "C20" = 0.7400577619174555
    "D20" = 0.02893492280009013
    "E20" = 0.6095070390605883
    "F20" = 1.679196586759787
    "I20" = 0.9226455484459137
    "B21" = 1.943823792032958
    "C21" = 0.8317703991930898
    "D21" = 0.02983526385366986
    "E21" = 0.6872091022820541
    "F21" = 1.792662048340901
    "I21" = 0.9664832725464976
    "B22" = 2.075488690337409
    "C22" = 0.8919249609699023
    "D22" = 0.03045274161284794
    "E22" = 0.738273492567771
    "F22" = 1.868012203187078
    "I22" = 0.9958104620327504
    "B23" = 2.005166434535226
    "C23" = 0.859798752631491
    "D23" = 0.03012045635393434
    "E23" = 0.7109922513340479
    "F23" = 1.827685029803433
    "I23" = 0.9800949631544853
    "B24" = 1.739913128867101
    "C24" = 0.7385621960838762
    "D24" = 0.02892068272055326
    "E24" = 0.6082414728318923
    "F24" = 1.677361735266004
    "I24" = 0.9219402219909512
    "B25" = 1.455966638430652
    "C25" = 0.6086220630136268
    "D25" = 0.02774455462341763
    "E25" = 0.4984734630056238
    "F25" = 1.520126681014943
    "I25" = 0.8620056120074793
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
